$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card2")

# Fix header text for column N (remove trailing space)
$ws.Range("N1").Value = "Correction"

# New header for column O - copy formatting from a neighboring header cell
# (bold, centered, bordered) then set its own text.
$ws.Range("N1").Copy($ws.Range("O1"))
$ws.Range("O1").Value = "Serviced by "

# Fill previously-empty N column data cells with "nan" (matches existing pattern in rows 2-13)
for ($r = 2; $r -le 13; $r++) {
    $ws.Cells.Item($r, 14).Value = "nan"
}

# Create the new (empty) O column data cells for rows 2-13, matching the
# default/unstyled look of the rest of the data rows.
for ($r = 2; $r -le 13; $r++) {
    $ws.Cells.Item($r, 15).Style = "Normal"
}
